# open ai bilan integratsiya qo'shildi
#
# Users sheet: bump Firdavs' salary, replace Nodira's row with Lola's
# details, and remove the last two rows (Bekzod, Shaxriyor).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Firdavs (row 2): salary 5000 -> 6000
$ws.Range("E2").Value = 6000

# Row 4: Nodira -> Lola (new email/department casing/salary)
$ws.Range("B4").Value = "Lola"
$ws.Range("C4").Value = "lola@gmail.com"
$ws.Range("D4").Value = "Hr"
$ws.Range("E4").Value = 4000

# Remove rows 5 (Bekzod) and 6 (Shaxriyor) entirely, shrinking the table
# from A1:E6 down to A1:E4.
$ws.Rows("5:6").Delete()
